$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.055.22"
$ws.Range("E2").Value = "  -4.07%  "
$ws.Range("D3").Value = "1.963.87"
$ws.Range("E3").Value = "  -6.11%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'327.56"
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.5001"
$ws.Range("E7").Value = "  -5.81%  "
$ws.Range("D8").Value = "'0.4214"
$ws.Range("E8").Value = "  -3.91%  "
$ws.Range("D9").Value = "'53.02"
$ws.Range("E9").Value = "  -2.50%  "
$ws.Range("D10").Value = "'0.09192"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").Value = "'1.099"
$ws.Range("E11").Value = "  -6.51%  "
$ws.Range("D12").Value = "'22.91"
$ws.Range("E12").Value = "  -7.53%  "
$ws.Range("D13").Value = "2.020.87"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "'7.866"
$ws.Range("E14").Value = "  -8.20%  "
$ws.Range("D15").Value = "'6.439"
$ws.Range("E15").Value = "  -6.52%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "'0.00001100"
$ws.Range("E17").Value = "  -5.01%  "
$ws.Range("D18").Value = "'91.41"
$ws.Range("E18").Value = "  -10.01%  "
$ws.Range("D19").Value = "'0.06678"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "'19.22"
$ws.Range("E20").Value = "  -9.09%  "
$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "'5.965"
$ws.Range("E22").Value = "  -6.10%  "
$ws.Range("D23").Value = "29.096.02"
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("D24").Value = "'12.09"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").Value = "'2.284"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "2.226.92"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'156.44"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D28").Value = "'20.61"
$ws.Range("E28").Value = "  -5.55%  "
$ws.Range("D29").Value = "'6.197"
$ws.Range("E29").Value = "  -9.83%  "
$ws.Range("D30").Value = "'2.261"
$ws.Range("E30").Value = "  -9.56%  "
$ws.Range("D31").Value = "'126.72"
$ws.Range("E31").Value = "  -5.18%  "
$ws.Range("D32").Value = "'1.042"
$ws.Range("E32").Value = "  -8.10%  "
$ws.Range("D33").Value = "'0.09852"
$ws.Range("E33").Value = "  -6.34%  "
$ws.Range("D34").Value = "'1.528"
$ws.Range("E34").Value = "  -8.45%  "
$ws.Range("D35").Value = "'5.777"
$ws.Range("E35").Value = "  -7.95%  "
$ws.Range("D36").Value = "'3.674"
$ws.Range("E36").Value = "  -6.00%  "
$ws.Range("D37").Value = "'0.02423"
$ws.Range("E37").Value = "  -7.66%  "
$ws.Range("D38").Value = "'1.300"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").Value = "'8.970"
$ws.Range("E39").Value = "  -11.18%  "
$ws.Range("D40").Value = "'0.06321"
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("D41").Value = "'0.6444"
$ws.Range("E41").Value = "  -7.37%  "
$ws.Range("D42").Value = "'11.44"
$ws.Range("E42").Value = "  -9.01%  "
$ws.Range("D43").Value = "'0.1990"
$ws.Range("E43").Value = "  -10.04%  "
$ws.Range("D44").Value = "'1.005"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "'0.6244"
$ws.Range("E45").Value = "  -7.73%  "
$ws.Range("D46").Value = "'13.30"
$ws.Range("E46").Value = "  -7.23%  "
$ws.Range("D47").Value = "'2.180"
$ws.Range("E47").Value = "  -8.82%  "
$ws.Range("D48").Value = "'1.291"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "'3.467"
$ws.Range("E49").Value = "  -4.50%  "
$ws.Range("E50").Value = "  -5.20%  "
$ws.Range("D51").Value = "'0.06931"
$ws.Range("E51").Value = "  -4.59%  "
